$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("N1_D40")
$ws.Range("E2").Value = 0.018
$ws.Range("F2").Value = 15.51
$ws.Range("E3").Value = 0.018
$ws.Range("F3").Value = 15.37
$ws.Range("F4").Value = 15.32
$ws.Range("F5").Value = 15.3
$ws.Range("E6").Value = 0.018
$ws.Range("F6").Value = 15.3
$ws.Range("F7").Value = 15.43
$ws.Range("F8").Value = 15.46
$ws.Range("F9").Value = 15.45
$ws.Range("E10").Value = 0.018
$ws.Range("F10").Value = 15.37
$ws.Range("E11").Value = 0.018
$ws.Range("F11").Value = 15.47
$ws.Range("E12").Value = 0.0178
$ws.Range("F12").Value = 15.398

$ws = $wb.Worksheets.Item("N1_D60")
$ws.Range("F2").Value = 14.87
$ws.Range("E3").Value = 0.027
$ws.Range("F3").Value = 14.76
$ws.Range("E4").Value = 0.026
$ws.Range("F4").Value = 14.97
$ws.Range("E5").Value = 0.027
$ws.Range("E6").Value = 0.027
$ws.Range("F6").Value = 14.83
$ws.Range("F7").Value = 14.76
$ws.Range("E8").Value = 0.026
$ws.Range("F8").Value = 14.87
$ws.Range("F9").Value = 14.75
$ws.Range("E10").Value = 0.027
$ws.Range("F10").Value = 14.9
$ws.Range("F11").Value = 14.76
$ws.Range("E12").Value = 0.0264
$ws.Range("F12").Value = 14.834

$ws = $wb.Worksheets.Item("N1_D80")
$ws.Range("F2").Value = 18.38
$ws.Range("F3").Value = 18.42
$ws.Range("E4").Value = 0.04
$ws.Range("F4").Value = 18.43
$ws.Range("E5").Value = 0.04
$ws.Range("F5").Value = 18.31
$ws.Range("F6").Value = 18.32
$ws.Range("E7").Value = 0.04
$ws.Range("F7").Value = 18.31
$ws.Range("E8").Value = 0.04
$ws.Range("F8").Value = 18.31
$ws.Range("E9").Value = 0.041
$ws.Range("F9").Value = 18.42
$ws.Range("E10").Value = 0.041
$ws.Range("F10").Value = 18.31
$ws.Range("E11").Value = 0.039
$ws.Range("F11").Value = 18.3
$ws.Range("E12").Value = 0.0401
$ws.Range("F12").Value = 18.351

$ws = $wb.Worksheets.Item("N1_D100")
$ws.Range("E2").Value = 0.058
$ws.Range("F2").Value = 16.47
$ws.Range("E3").Value = 0.06
$ws.Range("F3").Value = 16.5
$ws.Range("E4").Value = 0.06
$ws.Range("F4").Value = 16.57
$ws.Range("E5").Value = 0.057
$ws.Range("F5").Value = 16.5
$ws.Range("F6").Value = 16.43
$ws.Range("E7").Value = 0.057
$ws.Range("F7").Value = 16.42
$ws.Range("F8").Value = 16.57
$ws.Range("E9").Value = 0.059
$ws.Range("F9").Value = 16.51
$ws.Range("E10").Value = 0.059
$ws.Range("F10").Value = 16.57
$ws.Range("E11").Value = 0.061
$ws.Range("F11").Value = 16.48
$ws.Range("E12").Value = 0.0587
$ws.Range("F12").Value = 16.502
